# Apply "new variables for when teams submitted forecasts" edit:
# Adds two new reference rows (TournamentStart, ForecastisUpdate) describing
# two new survey variables, right after the existing last data row (103).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 104: TournamentStart -------------------------------------------------
$ws.Cells.Item(104, 1).Value = "TournamentStart"
$ws.Cells.Item(104, 2).Value = "Both"
$ws.Cells.Item(104, 3).Value = "Both"
$ws.Cells.Item(104, 4).Value = '"May", "October"'
$ws.Cells.Item(104, 4).WrapText = $true
$ws.Cells.Item(104, 5).Value = "Indicates whether the team started in May or October"

# --- Row 105: ForecastisUpdate ------------------------------------------------
$ws.Cells.Item(105, 1).Value = "ForecastisUpdate"
$ws.Cells.Item(105, 2).Value = "Academic"
$ws.Cells.Item(105, 3).Value = "Phase 2"
$ws.Cells.Item(105, 4).Value = "0 = no update, 1 = update"
$ws.Cells.Item(105, 4).WrapText = $true
$ws.Cells.Item(105, 5).Value = "Indicates whether the team's forecast for phase 2 is an update of their previous forecast (1) or a new forecast (0)"

# Reflect the author's final cursor position/selection in the saved view state.
[void]$ws.Range("C111").Select()
